$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly Mango price record was collected; insert it as a new row
# right above the current row 67, shifting all the existing rows below it
# down by one (the previous last row, 143, becomes 144).
$ws.Rows(67).Insert()

$ws.Cells.Item(67, 1).Value  = 5
$ws.Cells.Item(67, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(67, 3).Value  = "Maule"
$ws.Cells.Item(67, 4).Value  = 44799
$ws.Cells.Item(67, 5).Value  = 7
$ws.Cells.Item(67, 6).Value  = "Fruta"
$ws.Cells.Item(67, 7).Value  = 100108
$ws.Cells.Item(67, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(67, 9).Value  = 100108002
$ws.Cells.Item(67, 10).Value = "Mango"
$ws.Cells.Item(67, 11).Value = "Sin especificar"
$ws.Cells.Item(67, 12).Value = "Primera"
$ws.Cells.Item(67, 13).Value = 228
$ws.Cells.Item(67, 14).Value = 9000
$ws.Cells.Item(67, 15).Value = 10000
$ws.Cells.Item(67, 16).Value = 9439
$ws.Cells.Item(67, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(67, 18).Value = "México"
$ws.Cells.Item(67, 19).Value = 2360
$ws.Cells.Item(67, 20).Value = 4
